$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.995.77"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "2.037.19"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "228.71"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("D7").Value = "60.81"
$ws.Range("E7").Value = "  +3.47%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -0.79%  "
$ws.Range("E10").Value = "  +1.21%  "
$ws.Range("E11").Value = "  +0.82%  "
$ws.Range("D12").Value = "14.63"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").Value = "2.336.43"
$ws.Range("E13").Value = "  -0.78%  "
$ws.Range("D14").Value = "21.41"
$ws.Range("E14").Value = "  +2.73%  "
$ws.Range("D15").Value = "0.764"
$ws.Range("E15").Value = "  +1.77%  "
$ws.Range("D16").Value = "5.17"
$ws.Range("E16").Value = "  -2.03%  "
$ws.Range("D17").Value = "2.024.56"
$ws.Range("E17").Value = "  -1.85%  "
$ws.Range("D18").Value = "37.874.57"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").Value = "69.94"
$ws.Range("E19").Value = "  +0.42%  "
$ws.Range("D20").Value = "5.97"
$ws.Range("E20").Value = "  -4.75%  "
$ws.Range("D21").Value = "0.0₃0828"
$ws.Range("E21").Value = "  -0.78%  "
$ws.Range("D22").Value = "224.79"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "2.42"
$ws.Range("E24").Value = "  -0.44%  "
$ws.Range("E25").Value = "  +1.18%  "
$ws.Range("D26").Value = "167.21"
$ws.Range("E26").Value = "  +0.51%  "
$ws.Range("D27").Value = "9.31"
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("D28").Value = "0.129"
$ws.Range("E28").Value = "  -3.80%  "
$ws.Range("D29").Value = "18.95"
$ws.Range("E29").Value = "  -0.36%  "
$ws.Range("D30").Value = "1.29"
$ws.Range("E30").Value = "  -3.05%  "
$ws.Range("E31").Value = "  +1.05%  "
$ws.Range("D32").Value = "2.15"
$ws.Range("E32").Value = "  +4.70%  "
$ws.Range("D33").Value = "4.42"
$ws.Range("E33").Value = "  -2.13%  "
$ws.Range("D34").Value = "0.0608"
$ws.Range("E34").Value = "  -0.35%  "
$ws.Range("D35").Value = "4.52"
$ws.Range("E35").Value = "  -1.30%  "
$ws.Range("D36").Value = "6.49"
$ws.Range("E36").Value = "  +6.92%  "
$ws.Range("D37").Value = "2.29"
$ws.Range("E37").Value = "  -1.68%  "
$ws.Range("D38").Value = "3.27"
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.26%  "
$ws.Range("D40").Value = "1.527.05"
$ws.Range("E40").Value = "  +2.61%  "
$ws.Range("E41").Value = "  +4.28%  "
$ws.Range("E42").Value = "  +0.29%  "
$ws.Range("D43").Value = "96.41"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("D44").Value = "2.83"
$ws.Range("E44").Value = "  -1.67%  "
$ws.Range("D45").Value = "0.0916"
$ws.Range("E45").Value = "  -0.83%  "
$ws.Range("E46").Value = "  -1.59%  "
$ws.Range("D47").Value = "3.97"
$ws.Range("E47").Value = "  -3.79%  "
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").Value = "7.11"
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("D51").Value = "2.224.01"
$ws.Range("E51").Value = "  -0.85%  "
